# Applies the Fri Aug  9 17:28:01 UTC 2024 GitHub Actions cryptos-list refresh:
# updates the Price (D) / Volume(1h) (E) columns for each ranked coin row,
# re-sorts the Stacks/Bittensor pair (rows 42-43 swap rank), and replaces the
# WhiteBITCoin row (51) with the newly-ranked Maker entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Columns B/C/E never look like bare numbers, but D frequently does
    # ("506.31", "0.0000140", ...). Force text format first so Excel
    # keeps the literal string instead of silently reparsing it as a number.
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 2
$ws.Range("D2").Value = "60.354.94"
$ws.Range("E2").Value = "  +2.06%  "

# Row 3
$ws.Range("D3").Value = "2.587.60"
$ws.Range("E3").Value = "  +0.01%  "

# Row 4
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
Set-TextValue $ws.Range("D5") "506.31"
$ws.Range("E5").Value = "  -0.11%  "

# Row 6
Set-TextValue $ws.Range("D6") "153.13"
$ws.Range("E6").Value = "  -3.24%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.999"
$ws.Range("E7").Value = "  +0.50%  "

# Row 8
$ws.Range("E8").Value = "  -6.39%  "

# Row 9
$ws.Range("D9").Value = "2.597.26"
$ws.Range("E9").Value = "  +0.48%  "

# Row 10
Set-TextValue $ws.Range("D10") "6.64"
$ws.Range("E10").Value = "  +7.26%  "

# Row 11
$ws.Range("E11").Value = "  -0.37%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.345"
$ws.Range("E12").Value = "  +0.68%  "

# Row 13
$ws.Range("E13").Value = "  +0.96%  "

# Row 14
$ws.Range("D14").Value = "3.044.38"
$ws.Range("E14").Value = "  +1.89%  "

# Row 15
$ws.Range("D15").Value = "60.309.61"
$ws.Range("E15").Value = "  +2.01%  "

# Row 16
Set-TextValue $ws.Range("D16") "21.56"
$ws.Range("E16").Value = "  -1.94%  "

# Row 17
Set-TextValue $ws.Range("D17") "0.0000140"
$ws.Range("E17").Value = "  +1.20%  "

# Row 18
$ws.Range("D18").Value = "2.592.67"
$ws.Range("E18").Value = "  +0.84%  "

# Row 19
$ws.Range("E19").Value = "  +1.51%  "

# Row 20
Set-TextValue $ws.Range("D20") "345.85"
$ws.Range("E20").Value = "  +3.02%  "

# Row 21
Set-TextValue $ws.Range("D21") "10.34"
$ws.Range("E21").Value = "  -0.29%  "

# Row 22
Set-TextValue $ws.Range("D22") "6.11"
$ws.Range("E22").Value = "  +0.86%  "

# Row 23
$ws.Range("E23").Value = "  -0.80%  "

# Row 24
Set-TextValue $ws.Range("D24") "60.24"
$ws.Range("E24").Value = "  +0.04%  "

# Row 25
Set-TextValue $ws.Range("D25") "0.419"
$ws.Range("E25").Value = "  +0.67%  "

# Row 26
$ws.Range("E26").Value = "  -1.49%  "

# Row 27
$ws.Range("D27").Value = "2.698.97"
$ws.Range("E27").Value = "  +1.61%  "

# Row 28
$ws.Range("E28").Value = "  +0.40%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0844"
$ws.Range("E29").Value = "  +1.49%  "

# Row 30
$ws.Range("E30").Value = "  -0.60%  "

# Row 31
$ws.Range("E31").Value = "  +0.29%  "

# Row 32
Set-TextValue $ws.Range("D32") "19.29"
$ws.Range("E32").Value = "  -1.05%  "

# Row 33
Set-TextValue $ws.Range("D33") "154.60"
$ws.Range("E33").Value = "  -0.79%  "

# Row 34
$ws.Range("E34").Value = "  -0.96%  "

# Row 35
$ws.Range("E35").Value = "  +4.46%  "

# Row 36
Set-TextValue $ws.Range("D36") "3.99"
$ws.Range("E36").Value = "  +1.53%  "

# Row 37
Set-TextValue $ws.Range("D37") "1.19"
$ws.Range("E37").Value = "  -0.30%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.852"
$ws.Range("E38").Value = "  +18.51%  "

# Row 39
$ws.Range("E39").Value = "  -1.39%  "

# Row 40
$ws.Range("E40").Value = "  +0.12%  "

# Row 41
$ws.Range("E41").Value = "  +2.55%  "

# Row 42
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D42") "1.45"
$ws.Range("E42").Value = "  +1.03%  "

# Row 43
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D43") "296.70"
$ws.Range("E43").Value = "  +1.71%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.619"
$ws.Range("E44").Value = "  -1.37%  "

# Row 45
$ws.Range("E45").Value = "  -1.75%  "

# Row 46
Set-TextValue $ws.Range("D46") "0.0561"
$ws.Range("E46").Value = "  +0.21%  "

# Row 47
Set-TextValue $ws.Range("D47") "0.997"
$ws.Range("E47").Value = "  +0.23%  "

# Row 48
Set-TextValue $ws.Range("D48") "19.74"
$ws.Range("E48").Value = "  +0.95%  "

# Row 49
Set-TextValue $ws.Range("D49") "4.85"
$ws.Range("E49").Value = "  +0.28%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.0234"
$ws.Range("E50").Value = "  -1.92%  "

# Row 51
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.008.51"
$ws.Range("E51").Value = "  +1.69%  "
